# "express input corregidas hasta la 17"
# Corrects the "MEDICION" (column B) labels for a set of rows in Sheet1:
#  - med_01 (row 2): amplitude typo 2Vpp -> 6Vpp
#  - several rows: drop the erroneous "+ FR" suffix from the label
#    (these measurements were taken without the FR stage, unlike their
#    "+ FR" counterparts elsewhere in the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "FAA + LLAVE + FR - XA: Senoidal 500Hz 6Vpp - Duty 24%"
$ws.Range("B4").Value  = "FAA + LLAVE - XB: exp{-|t|} 10Hz 2V - Duty 20%"
$ws.Range("B5").Value  = "FAA + S&H - XB: exp{-|t|} 10Hz 2V - Duty 20%"
$ws.Range("B6").Value  = "FAA + S&H - XA: Senoidal 500Hz 2Vpp - Duty 71%"
$ws.Range("B7").Value  = "FAA + LLAVE - XA: Senoidal 500Hz 2Vpp - Duty 88%"
$ws.Range("B8").Value  = "FAA + LLAVE - XB: exp{-|t|} 10Hz 2V - Duty 30%"
$ws.Range("B10").Value = "FAA + S&H - XA: Senoidal 1KHz 2Vpp - Duty 14%"
$ws.Range("B11").Value = "FAA + S&H - XA: Senoidal 500Hz 2Vpp - Duty 24%"
$ws.Range("B12").Value = "FAA + LLAVE - XA: Senoidal 500Hz 2Vpp - Duty 24%"
$ws.Range("B14").Value = "FAA + LLAVE - XA: 3/2 Seno 500Hz 2Vpp - Duty 44%"
$ws.Range("B16").Value = "FAA + S&H - XA: 3/2 Seno 500Hz 2Vpp - Duty 8%"
$ws.Range("B22").Value = "FAA + LLA + S&H - XC - Duty LLA 70% - Duty S&H 10% - fs = 3KHz"

# Move the active selection to row 19 (whole row), matching the cursor
# position left behind by the author after making these edits.
$ws.Range("A19:XFD19").Select()
